# Generate Report for Handoff
# Adds a new source file (af1f09fe-b970-4b01-9e3a-5bd87979b265.md) to the
# localization status report, while the previously tracked file's UUID is
# renamed from bc815363-6bb6-46ea-a089-fffb3deeab5c.md to
# 11997f48-c54a-44ad-a860-376efe8dc576.md, on all three worksheets
# (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Common identifiers used across sheets
# ---------------------------------------------------------------------
$mdUuid1 = "11997f48-c54a-44ad-a860-376efe8dc576"
$mdUuid2 = "af1f09fe-b970-4b01-9e3a-5bd87979b265"
$md1 = "$mdUuid1.md"
$md2 = "$mdUuid2.md"
$hash1 = "38f602dae8dbfae6cf0cb715f8b260a3c6def621"
$hash2 = "3026bd34376f1ce1385c5026ca1e38890f340150"

$xlf1zh = "$mdUuid1.$hash1.zh-cn.xlf"
$xlf2zh = "$mdUuid2.$hash2.zh-cn.xlf"
$xlf1de = "$mdUuid1.$hash1.de-de.xlf"
$xlf2de = "$mdUuid2.$hash2.de-de.xlf"

$urlMd1 = "https://github.com/OpenLocalizationTest/oltest/blob/1aa2f7c6f3c37944d06030e5fa5ae7b211c5e79c/e2e/$md1"
$urlMd2 = "https://github.com/OpenLocalizationTest/oltest/blob/1aa2f7c6f3c37944d06030e5fa5ae7b211c5e79c/e2e/$md2"
$urlConfig = "https://github.com/OpenLocalizationTest/oltest/blob/1aa2f7c6f3c37944d06030e5fa5ae7b211c5e79c/.localization-config"

$urlXlf1zh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e5f3c3dce4ddd7e78a52a0df0008249d4711708c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlf1zh"
$urlXlf2zh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e5f3c3dce4ddd7e78a52a0df0008249d4711708c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlf2zh"
$urlXlf1de = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/82fb143824b42fc04e3da60df8d9cf2a448818f2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlf1de"
$urlXlf2de = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/82fb143824b42fc04e3da60df8d9cf2a448818f2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlf2de"

$readyForHandoff = "Ready for handoff"
$notLocalized = "Not to be localized"
$include = "Include"
$ignored = "Ignored"
$epoch = "0001-01-01 00:00:00"
$dateZh = "2016-03-08 06:33:12"
$dateDe = "2016-03-08 06:33:15"

# ---------------------------------------------------------------------
# Sheet 1 : "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Hyperlinks.Delete()

$ws1.Range("A2").Value = $md1
$ws1.Range("B2").Value = $readyForHandoff
$ws1.Range("C2").Value = $readyForHandoff

$ws1.Range("A3").Value = $md2
$ws1.Range("B3").Value = $readyForHandoff
$ws1.Range("C3").Value = $readyForHandoff

$ws1.Range("A4").Value = ".localization-config"
$ws1.Range("B4").Value = $notLocalized
$ws1.Range("C4").Value = $notLocalized

$ws1.Hyperlinks.Add($ws1.Range("A2"), $urlMd1, "", "", $md1) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), $urlMd2, "", "", $md2) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), $urlConfig, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2 : "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()

$ws2.Range("A2").Value = $md1
$ws2.Range("B2").Value = $readyForHandoff
$ws2.Range("C2").Value = $xlf1zh
$ws2.Range("D2").Value = $dateZh
$ws2.Range("G2").Value = $epoch
$ws2.Range("H2").Value = $include

$ws2.Range("A3").Value = $md2
$ws2.Range("B3").Value = $readyForHandoff
$ws2.Range("C3").Value = $xlf2zh
$ws2.Range("D3").Value = $dateZh
$ws2.Range("G3").Value = $epoch
$ws2.Range("H3").Value = $include

$ws2.Range("A4").Value = ".localization-config"
$ws2.Range("B4").Value = $notLocalized
$ws2.Range("D4").Value = $epoch
$ws2.Range("G4").Value = $epoch
$ws2.Range("H4").Value = $ignored

$ws2.Hyperlinks.Add($ws2.Range("A2"), $urlMd1, "", "", $md1) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), $urlXlf1zh, "", "", $xlf1zh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), $urlMd2, "", "", $md2) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), $urlXlf2zh, "", "", $xlf2zh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), $urlConfig, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet 3 : "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()

$ws3.Range("A2").Value = $md1
$ws3.Range("B2").Value = $readyForHandoff
$ws3.Range("C2").Value = $xlf1de
$ws3.Range("D2").Value = $dateDe
$ws3.Range("G2").Value = $epoch
$ws3.Range("H2").Value = $include

$ws3.Range("A3").Value = $md2
$ws3.Range("B3").Value = $readyForHandoff
$ws3.Range("C3").Value = $xlf2de
$ws3.Range("D3").Value = $dateDe
$ws3.Range("G3").Value = $epoch
$ws3.Range("H3").Value = $include

$ws3.Range("A4").Value = ".localization-config"
$ws3.Range("B4").Value = $notLocalized
$ws3.Range("D4").Value = $epoch
$ws3.Range("G4").Value = $epoch
$ws3.Range("H4").Value = $ignored

$ws3.Hyperlinks.Add($ws3.Range("A2"), $urlMd1, "", "", $md1) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), $urlXlf1de, "", "", $xlf1de) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), $urlMd2, "", "", $md2) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), $urlXlf2de, "", "", $xlf2de) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), $urlConfig, "", "", ".localization-config") | Out-Null
